# Updates cryptos list — prices (column D) and 1h volume % (column E)
# for rows 2-51. Values are forced to text via a leading apostrophe so
# numeric-looking strings (e.g. "583.76") stay text instead of becoming
# Excel numbers, matching the source data's inline-string cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.266.63"
$ws.Range("E2").Value = "'  -0.83%  "
$ws.Range("D3").Value = "'2.604.07"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'583.76"
$ws.Range("E5").Value = "'  +2.05%  "
$ws.Range("D6").Value = "'142.80"
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("E8").Value = "'  -0.90%  "
$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("E10").Value = "'  -1.91%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("E11").Value = "'  +0.01%  "
$ws.Range("E12").Value = "'  +1.06%  "
$ws.Range("D13").Value = "'3.064.02"
$ws.Range("E13").Value = "'  -0.34%  "
$ws.Range("D14").Value = "'24.61"
$ws.Range("E14").Value = "'  +4.28%  "
$ws.Range("D15").Value = "'60.253.33"
$ws.Range("E15").Value = "'  -0.89%  "
$ws.Range("E16").Value = "'  +0.03%  "
$ws.Range("D17").Value = "'2.606.87"
$ws.Range("E17").Value = "'  -0.74%  "
$ws.Range("E18").Value = "'  +1.08%  "
$ws.Range("E19").Value = "'  -1.74%  "
$ws.Range("D20").Value = "'347.52"
$ws.Range("E20").Value = "'  -0.25%  "
$ws.Range("E21").Value = "'  -2.20%  "
$ws.Range("E22").Value = "'  -0.30%  "
$ws.Range("D23").Value = "'0.537"
$ws.Range("E23").Value = "'  +4.00%  "
$ws.Range("D24").Value = "'63.80"
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "'  +0.48%  "
$ws.Range("E26").Value = "'  -0.71%  "
$ws.Range("E27").Value = "'  +2.92%  "
$ws.Range("E28").Value = "'  +1.19%  "
$ws.Range("D29").Value = "'0.0₃0795"
$ws.Range("E29").Value = "'  -0.13%  "
$ws.Range("E30").Value = "'  +4.72%  "
$ws.Range("E31").Value = "'  +1.59%  "
$ws.Range("E32").Value = "'  +0.11%  "
$ws.Range("E33").Value = "'  -0.57%  "
$ws.Range("D34").Value = "'1.31"
$ws.Range("E34").Value = "'  +9.57%  "
$ws.Range("D35").Value = "'4.24"
$ws.Range("E35").Value = "'  -0.68%  "
$ws.Range("E36").Value = "'  +3.27%  "
$ws.Range("E37").Value = "'  +2.64%  "
$ws.Range("D38").Value = "'316.94"
$ws.Range("E38").Value = "'  +6.65%  "
$ws.Range("D39").Value = "'38.25"
$ws.Range("E39").Value = "'  +1.49%  "
$ws.Range("E40").Value = "'  +1.95%  "
$ws.Range("E41").Value = "'  -0.85%  "
$ws.Range("D42").Value = "'135.81"
$ws.Range("E42").Value = "'  -2.12%  "
$ws.Range("E43").Value = "'  +0.93%  "
$ws.Range("E44").Value = "'  +0.20%  "
$ws.Range("D45").Value = "'19.92"
$ws.Range("E45").Value = "'  +1.34%  "
$ws.Range("E46").Value = "'  +0.11%  "
$ws.Range("D47").Value = "'0.0547"
$ws.Range("E47").Value = "'  -0.86%  "
$ws.Range("D48").Value = "'4.96"
$ws.Range("E48").Value = "'  +3.05%  "
$ws.Range("E49").Value = "'  -0.06%  "
$ws.Range("D50").Value = "'19.94"
$ws.Range("E50").Value = "'  +1.74%  "
$ws.Range("E51").Value = "'  +0.31%  "
